$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ForgotPassword_Tests")

# --- Row 2: add Actual Result for FP-001 (new column H) ---
$ws.Range("H2").Value = "Email input field is present on the forgot password page."
$ws.Range("B2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# --- Row 3: cosmetic font cleanup on F3 (content/value unchanged) ---
$ws.Range("B2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Rows.Item(3).RowHeight = 145.8

# --- Row 4: refine Expected Result wording + add Actual Result for FP-003 ---
$ws.Range("G4").Value = "An error message  should be displayed  for invalid email format."
$ws.Range("H4").Value = "An error message is displayed for invalid email format as expected."
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

# --- Row 5: brand-new FP-004 test case (Verify error message for empty email field) ---
$ws.Range("A5").Value = "FP-004"
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("B5").Value = "Verify error message for empty email field"
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("C5").Value = "To ensure that an error message is displayed if the email field is left empty."
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("D5").Value = "Open the Url - https://magento.softwaretestingboard.com/."
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$steps5 = @"
1. Navigate to 'sign in' page
2.Navigate to 'Forgot Password' page.
3.Don't enter anything in email address field. 
3. Click the submit button.
4. Verify the displayed error message.
"@
$ws.Range("E5").Value = $steps5
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)

$data5 = @"
{
        "Email": ""
    }
"@
$ws.Range("F5").Value = $data5
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("G5").Value = "An error message  should be displayed."
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("H5").Value = "An error message is displayed for empty email address field as expected."
$ws.Range("G4").Copy()
$ws.Range("H5").PasteSpecial(-4122)

$ws.Rows.Item(5).RowHeight = 97.2

# --- Selection marker matches authored workbook ---
$ws.Activate()
$ws.Range("D3").Select()
